$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45966
$ws.Range("B2").Value = 30.89
$ws.Range("C2").Value = 24.59
$ws.Range("D2").Value = 17.21
$ws.Range("E2").Value = 13.87
$ws.Range("F2").Value = 7
$ws.Range("G2").Value = 13.89
$ws.Range("H2").Value = 27.58
$ws.Range("I2").Value = 33.99
$ws.Range("J2").Value = 34.07
$ws.Range("K2").Value = 21.28
$ws.Range("L2").Value = 8.56
$ws.Range("M2").Value = 3.72
$ws.Range("N2").Value = 5.32
$ws.Range("O2").Value = 5.26
$ws.Range("P2").Value = 9.76
$ws.Range("Q2").Value = 16.72
$ws.Range("R2").Value = 26.66
$ws.Range("S2").Value = 33.41
$ws.Range("T2").Value = 46.22
$ws.Range("U2").Value = 65.88
$ws.Range("V2").Value = 70.5
$ws.Range("W2").Value = 66.8
$ws.Range("X2").Value = 47.94
$ws.Range("Y2").Value = 34.98
$ws.Range("Z2").Value = 27.75
$ws.Range("AB2").Value = 55.06
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 68.65000000000001
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 56.05
$ws.Range("AG2").Value = "1h-16h"
